$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") from 45175 -> 45177 for all data rows (2..494)
for ($r = 2; $r -le 494; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45177
    }
}

# 2) Row 494 gains an explicit custom row height (ht="15" customHeight="1")
$ws.Rows.Item(494).RowHeight = 15

# 3) Append new row 495 with the new entry
$ws.Range("A495").Value = "A 41757-2023"
$ws.Range("B495").Value = 45176
$ws.Range("C495").Value = 45177
$ws.Range("D495").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E495").Value = "GISLAVED"
$ws.Range("F495").Value = "Sveaskog"
$ws.Range("G495").Value = 2.3
$ws.Range("H495").Value = 0
$ws.Range("I495").Value = 0
$ws.Range("J495").Value = 0
$ws.Range("K495").Value = 0
$ws.Range("L495").Value = 0
$ws.Range("M495").Value = 0
$ws.Range("N495").Value = 0
$ws.Range("O495").Value = 0
$ws.Range("P495").Value = 0
$ws.Range("Q495").Value = 0
$ws.Range("R495").Value = ""

# Match formatting of the row above (date format on B/C, wrap-text on R)
$ws.Range("B495").NumberFormat = $ws.Range("B494").NumberFormat
$ws.Range("C495").NumberFormat = $ws.Range("C494").NumberFormat
$ws.Range("R495").WrapText = $ws.Range("R494").WrapText
